$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

$promptText = @"
 Given is the adjacency matrix for a weighted undirected graph containing 20 nodes labelled A to T. The value corresponding to each row M and column N represents the cost of travelling between the two nodes, where 0 means no connection.   

Consider some examples

Example 1: what is the least cost path from node A to node H?
   A B C D E F G H I J K L
 A 0 3 0 0 0 0 0 0 0 0 0 0
 B 3 0 4 0 0 0 0 0 3 2 3 3
 C 0 4 0 2 0 0 0 2 0 0 0 0
 D 0 0 2 0 1 1 0 0 0 0 0 0
 E 0 0 0 1 0 0 0 0 0 0 0 0
 F 0 0 0 1 0 0 3 0 0 0 0 0
 G 0 0 0 0 0 3 0 0 0 0 0 0
 H 0 0 2 0 0 0 0 0 0 0 0 0
 I 0 3 0 0 0 0 0 0 0 3 0 0
 J 0 2 0 0 0 0 0 0 3 0 0 0
 K 0 3 0 0 0 0 0 0 0 0 0 3
 L 0 3 0 0 0 0 0 0 0 0 3 0

Solution: A -> B -> C -> H
        

Example 2: what is the least cost path from node A to node J?
   A B C D E F G H I J K L
 A 0 3 4 0 0 2 0 0 0 0 0 0
 B 3 0 0 0 0 0 0 0 0 0 0 0
 C 4 0 0 2 0 0 0 0 0 0 0 0
 D 0 0 2 0 5 0 0 0 0 0 0 0
 E 0 0 0 5 0 0 0 0 0 0 0 0
 F 2 0 0 0 0 0 4 0 0 0 0 1
 G 0 0 0 0 0 4 0 5 0 2 3 0
 H 0 0 0 0 0 0 5 0 2 0 0 0
 I 0 0 0 0 0 0 0 2 0 0 0 0
 J 0 0 0 0 0 0 2 0 0 0 0 0
 K 0 0 0 0 0 0 3 0 0 0 0 3
 L 0 0 0 0 0 1 0 0 0 0 3 0

Solution: A -> F -> G -> J
        

Example 3: what is the least cost path from node A to node F?
   A B C D E F G H I J K
 A 0 4 0 5 0 0 3 0 0 0 0
 B 4 0 5 0 0 0 0 0 0 0 0
 C 0 5 0 0 0 0 0 0 0 0 0
 D 5 0 0 0 2 0 0 0 0 0 0
 E 0 0 0 2 0 5 0 4 0 0 0
 F 0 0 0 0 5 0 0 0 0 0 5
 G 3 0 0 0 0 0 0 0 0 0 0
 H 0 0 0 0 4 0 0 0 3 0 0
 I 0 0 0 0 0 0 0 3 0 3 0
 J 0 0 0 0 0 0 0 0 3 0 1
 K 0 0 0 0 0 5 0 0 0 1 0

Solution: A -> D -> E -> F
        
 Given these examples, answer the following quesiton.

what is the least cost path from node A to node J?

   A B C D E F G H I J K L M N O P Q R S T
 A 0 5 0 0 0 0 0 0 0 0 0 0 0 4 0 0 0 0 0 0
 B 5 0 2 0 2 0 0 0 0 0 4 0 0 0 0 0 0 0 0 0
 C 0 2 0 5 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 D 0 0 5 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 E 0 2 0 0 0 3 1 0 0 0 0 0 0 0 0 0 0 0 0 0
 F 0 0 0 0 3 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 G 0 0 0 0 1 0 0 2 5 0 0 0 0 0 0 0 0 0 0 1
 H 0 0 0 0 0 0 2 0 0 0 0 0 0 0 0 0 0 0 0 0
 I 0 0 0 0 0 0 5 0 0 4 0 0 0 0 0 0 0 0 0 0
 J 0 0 0 0 0 0 0 0 4 0 0 0 0 0 3 0 0 0 0 0
 K 0 4 0 0 0 0 0 0 0 0 0 3 0 0 0 0 0 0 0 0
 L 0 0 0 0 0 0 0 0 0 0 3 0 4 0 0 0 0 0 0 0
 M 0 0 0 0 0 0 0 0 0 0 0 4 0 3 0 0 0 0 0 0
 N 4 0 0 0 0 0 0 0 0 0 0 0 3 0 0 0 0 0 0 0
 O 0 0 0 0 0 0 0 0 0 3 0 0 0 0 0 2 0 0 0 0
 P 0 0 0 0 0 0 0 0 0 0 0 0 0 0 2 0 2 0 0 0
 Q 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 2 0 1 0 0
 R 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0 1 0
 S 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0 1
 T 0 0 0 0 0 0 1 0 0 0 0 0 0 0 0 0 0 0 1 0
    
"@

# --- Add new sheets "o_20" and "o_20_jumbled" after the existing "o_10" sheet ---
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "o_20"
$ws3 = $wb.Worksheets.Add($null, $ws2)
$ws3.Name = "o_20_jumbled"

# --- sheet1 ("o_10"): add column E header + row2 value ---
$ws1.Range("E1").Value = "evaluator_partial_correctness"
$ws1.Range("D1").Copy() | Out-Null
$ws1.Range("E1").PasteSpecial(-4122) | Out-Null

$ws1.Range("A2").Value = $promptText
$ws1.Range("B2").Value = "A -> B -> E -> G -> I -> J"
$ws1.Range("C2").Value = "From the adjacency matrix, we can determine the least cost path from node A to node J by finding the shortest path using an algorithm such as Dijkstra's algorithm."
$ws1.Range("D2").Value = "Correct"
$ws1.Range("E2").Value = "Output: 0/6"
$ws1.Range("A2").EntireRow.AutoFit()

# --- sheet2 ("o_20"): header row only, same style as sheet1 ---
$ws2.Range("A1").Value = "prompt"
$ws2.Range("B1").Value = "solution"
$ws2.Range("C1").Value = "llm_response"
$ws2.Range("D1").Value = "evaluator_response"
$ws2.Range("E1").Value = "evaluator_partial_correctness"
$ws1.Range("A1:D1").Copy() | Out-Null
$ws2.Range("A1:D1").PasteSpecial(-4122) | Out-Null
$ws1.Range("D1").Copy() | Out-Null
$ws2.Range("E1").PasteSpecial(-4122) | Out-Null

# --- sheet3 ("o_20_jumbled"): header row only, same style as sheet1 ---
$ws3.Range("A1").Value = "prompt"
$ws3.Range("B1").Value = "solution"
$ws3.Range("C1").Value = "llm_response"
$ws3.Range("D1").Value = "evaluator_response"
$ws3.Range("E1").Value = "evaluator_partial_correctness"
$ws1.Range("A1:D1").Copy() | Out-Null
$ws3.Range("A1:D1").PasteSpecial(-4122) | Out-Null
$ws1.Range("D1").Copy() | Out-Null
$ws3.Range("E1").PasteSpecial(-4122) | Out-Null

# restore original active sheet/tab selection
$ws1.Activate()

Write-Output "done"
